$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new "season record" columns after the existing data
# (Unnamed: 28 is in AC). Copy the header style from AC1 so the new
# headers (Wins/Losses/Ties) match the existing bold/bordered header look.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate every data row (2-50) with the team's season record.
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 71
    $ws.Cells.Item($r, 31).Value = 91
    $ws.Cells.Item($r, 32).Value = 0
}
